# Fix issue with chosen trenching type definition.
#
# - Units sheet: remove the now-redundant unit rows for
#   device.control_signal_type, options.user_installation_tool and
#   project.selected_installation_tool (these are categorical fields with
#   no physical unit; the valid-value list for the *selected* trenching
#   tool was missing instead).
# - Valid Values sheet: add the missing valid-value row for
#   project.selected_installation_tool (the chosen trenching type),
#   mirroring the existing options.user_installation_tool row.
# - Selection/view bookkeeping: Units becomes the active sheet/tab instead
#   of Tables, and ROOT/Units/Valid Values selections move to reflect the
#   edit location.

$wb = $excel.ActiveWorkbook

$units = $wb.Worksheets.Item("Units")
$root = $wb.Worksheets.Item("ROOT")
$tables = $wb.Worksheets.Item("Tables")
$validValues = $wb.Worksheets.Item("Valid Values")

# --- Units sheet: delete the three obsolete rows (bottom-up so row
#     numbers of the earlier deletions stay valid) ---
$units.Rows(37).Delete()
$units.Rows(30).Delete()
$units.Rows(16).Delete()

# --- Valid Values sheet: append the valid values for the chosen
#     (selected) installation/trenching tool ---
$validValues.Range("A7").Value = "project.selected_installation_tool"
$validValues.Range("B7").Value = "Jetting"
$validValues.Range("C7").Value = "Ploughing"
$validValues.Range("D7").Value = "Cutting"
$validValues.Range("E7").Value = "Dredging"
$validValues.Range("B7:E7").Style = $validValues.Range("B5:E5").Style
$validValues.Range("A7").Select()

# --- View / selection bookkeeping ---
$units.Range("A16:XFD16").Select()
$root.Range("A52").Select()

# Units becomes the active sheet/tab (was Tables)
$units.Activate()
